$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "r321"
$ws.Range("B9").Value = "barry"
$ws.Range("C9").Value = "testing the changes"
$ws.Range("D9").Value = "2025-09-30 17:11:53"
